$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Data values: Madagascar/French hierarchy -> Sierra Leone/English hierarchy
# ---------------------------------------------------------------------------
$ws.Range("A2:A6").Value = "eng"

$ws.Range("C3").Value = "Region"
$ws.Range("C4").Value = "District"
$ws.Range("C5").Value = "Chiefdom"
$ws.Range("C6").Value = "Section"

# New row 7 (Village) - copy formatting from row 6 first, then set values
$ws.Range("A6:D6").Copy() | Out-Null
$ws.Range("A7:D7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A7").Value = "eng"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "Village"
$ws.Range("D7").Value = $true
